$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A213").Value = "IMX-USD"
$ws.Range("A214").Value = "TAO-USD"
$ws.Range("A215").Value = "GRT-USD"
